$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns stay plain text so
# dotted/percent strings are not reinterpreted as numbers/dates.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "93.999.48"
$ws.Range("E2").Value = "  +1.55%  "

# Row 3
$ws.Range("D3").Value = "3.091.50"
$ws.Range("E3").Value = "  -0.66%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "234.54"
$ws.Range("E5").Value = "  -3.11%  "

# Row 6
$ws.Range("D6").Value = "609.73"
$ws.Range("E6").Value = "  -0.81%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "0.380"
$ws.Range("E8").Value = "  -5.01%  "

# Row 9
$ws.Range("E9").Value = "  -0.06%  "

# Row 10
$ws.Range("D10").Value = "0.820"
$ws.Range("E10").Value = "  +12.01%  "

# Row 11
$ws.Range("D11").Value = "3.086.65"
$ws.Range("E11").Value = "  -0.75%  "

# Row 12
$ws.Range("E12").Value = "  -3.49%  "

# Row 13
$ws.Range("D13").Value = "93.827.26"
$ws.Range("E13").Value = "  +1.65%  "

# Row 14
$ws.Range("D14").Value = "0.0000240"
$ws.Range("E14").Value = "  -6.12%  "

# Row 15
$ws.Range("D15").Value = "34.13"
$ws.Range("E15").Value = "  -0.84%  "

# Row 16
$ws.Range("B16").Value = "Toncoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D16").Value = "5.24"
$ws.Range("E16").Value = "  -4.77%  "

# Row 17
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.665.74"
$ws.Range("E17").Value = "  -0.67%  "

# Row 18
$ws.Range("D18").Value = "3.128.80"
$ws.Range("E18").Value = "  +1.21%  "

# Row 19
$ws.Range("E19").Value = "  -1.01%  "

# Row 20
$ws.Range("D20").Value = "14.67"
$ws.Range("E20").Value = "  -0.52%  "

# Row 21
$ws.Range("D21").Value = "5.79"
$ws.Range("E21").Value = "  -0.17%  "

# Row 22
$ws.Range("D22").Value = "442.23"
$ws.Range("E22").Value = "  -1.24%  "

# Row 23
$ws.Range("E23").Value = "  -6.27%  "

# Row 24
$ws.Range("E24").Value = "  -5.81%  "

# Row 25
$ws.Range("D25").Value = "8.26"
$ws.Range("E25").Value = "  +4.66%  "

# Row 26
$ws.Range("D26").Value = "5.53"
$ws.Range("E26").Value = "  -4.12%  "

# Row 27
$ws.Range("D27").Value = "84.75"
$ws.Range("E27").Value = "  -2.77%  "

# Row 28
$ws.Range("D28").Value = "11.95"
$ws.Range("E28").Value = "  +1.71%  "

# Row 29
$ws.Range("D29").Value = "3.255.64"
$ws.Range("E29").Value = "  -0.61%  "

# Row 30
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.16%  "

# Row 31
$ws.Range("D31").Value = "0.254"
$ws.Range("E31").Value = "  +9.58%  "

# Row 32
$ws.Range("D32").Value = "0.180"
$ws.Range("E32").Value = "  +6.48%  "

# Row 33
$ws.Range("E33").Value = "  -10.54%  "

# Row 34
$ws.Range("D34").Value = "9.27"
$ws.Range("E34").Value = "  -0.29%  "

# Row 35
$ws.Range("D35").Value = "0.995"
$ws.Range("E35").Value = "  -0.34%  "

# Row 36
$ws.Range("D36").Value = "7.80"
$ws.Range("E36").Value = "  -3.39%  "

# Row 37
$ws.Range("D37").Value = "0.159"
$ws.Range("E37").Value = "  -4.26%  "

# Row 38
$ws.Range("D38").Value = "25.61"
$ws.Range("E38").Value = "  -2.08%  "

# Row 39
$ws.Range("E39").Value = "  -1.96%  "

# Row 40
$ws.Range("D40").Value = "0.445"
$ws.Range("E40").Value = "  +0.67%  "

# Row 41
$ws.Range("D41").Value = "23.93"
$ws.Range("E41").Value = "  +3.75%  "

# Row 42
$ws.Range("E42").Value = "  -2.65%  "

# Row 43
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "467.24"
$ws.Range("E43").Value = "  -3.09%  "

# Row 44
$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D44").Value = "3.68"
$ws.Range("E44").Value = "  -13.56%  "

# Row 45
$ws.Range("E45").Value = "  -0.02%  "

# Row 46
$ws.Range("D46").Value = "3.12"
$ws.Range("E46").Value = "  -10.99%  "

# Row 47
$ws.Range("D47").Value = "161.51"

# Row 48
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.86"
$ws.Range("E48").Value = "  -2.86%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.680"
$ws.Range("E49").Value = "  -1.67%  "

# Row 50
$ws.Range("D50").Value = "43.74"
$ws.Range("E50").Value = "  -0.92%  "

# Row 51
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  +0.00%  "
